# Add a new worksheet "2017-05-04" (the latest daily diagnostics snapshot)
# positioned after the existing sheets, and populate it with the Data Entry
# Performance breakdown, mirroring the structure of the prior daily sheets.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing worksheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2017-05-04"

$data = @(
    @("CreatedBy", "DataCompleted", "DataIncomplete", "Total"),
    @("Dr Faeiz", 1, 0, 1),
    @("delina", 2, 0, 2),
    @("DR ANG", 1, 1, 2),
    @("Default Administrator", 1, 2, 3),
    @("Aziani", 9, 0, 9),
    @("Dr Saravanan", 10, 0, 10),
    @("Faeiz", 12, 0, 12),
    @("Dr Masliyana", 21, 0, 21),
    @("Hui Che", 21, 0, 21),
    @("Dr Faeiz Syezri Adzmin bin Jaaffar", 23, 0, 23),
    @("Ang SH", 28, 0, 28),
    @("Aimi Nadiah Jamel", 27, 2, 29),
    @("Fadzli", 30, 0, 30),
    @("Afiq Firdaus", 30, 1, 31),
    @("Jennifer Kaur", 31, 0, 31),
    @("Sunita", 32, 0, 32),
    @("Suhayl", 35, 0, 35),
    @("Ling Kuok Wei", 38, 0, 38),
    @("Eliza", 41, 0, 41),
    @("Adlan", 41, 1, 42),
    @("Musfirah", 43, 0, 43),
    @("Saravanan", 43, 0, 43),
    @("Noor hidayah", 45, 0, 45),
    @("Annas", 46, 0, 46),
    @("Delina", 51, 0, 51),
    @("Siti Aminah", 51, 0, 51),
    @("Aisyah", 52, 1, 53),
    @("Aizat", 53, 1, 54),
    @("Nursyuhaida", 54, 1, 55),
    @("Munirah", 58, 0, 58),
    @("Yhyviyaa", 60, 0, 60),
    @("Hui Yi", 61, 1, 62),
    @("Philip", 63, 0, 63),
    @("Noor Amalina", 64, 0, 64),
    @("Hadi", 66, 0, 66),
    @("Izzat", 67, 0, 67),
    @("Helmi", 68, 0, 68),
    @("Danial", 79, 0, 79),
    @("Natrah", 79, 0, 79),
    @("Syahirah", 81, 0, 81),
    @("Izzati", 82, 0, 82),
    @("Amira", 85, 0, 85),
    @("Michelle", 88, 0, 88),
    @("Tan khai shin", 88, 0, 88),
    @("Dr Richard", 89, 0, 89),
    @("Hooi Fan", 90, 0, 90),
    @("Yvonne", 89, 1, 90),
    @("Thivashini", 92, 0, 92),
    @("Mardhiah", 94, 0, 94),
    @("Pui Yee", 95, 0, 95),
    @("Aminiril anisah", 98, 0, 98),
    @("Geetha Krishnan", 103, 0, 103),
    @("Sree Durga", 104, 0, 104),
    @("Jocelyn", 109, 0, 109),
    @("Nabilah Iffah", 111, 0, 111),
    @("Hoong Ping", 113, 0, 113),
    @("Aishah", 114, 0, 114),
    @("Jia yi", 128, 0, 128),
    @("Nurjannah", 132, 0, 132),
    @("Yi Shin", 170, 0, 170),
    @("Jacelyn", 230, 0, 230),
    @("Total", 3922, 12, 3934)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 1
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
